# Update countries & provincias Spain
#
# Updates the COVID country-stats sheet:
#  - refreshes the "last updated" timestamp
#  - refreshes several countries' case/death counters
#  - Angola's case count overtakes Mali's -> the two swap rank (row 134/135)
#  - "Islas Malvinas" overtakes Montserrat -> the two swap rank (row 214/215)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 22:49"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 6374204
$ws.Range("C4").Value = 38960
$ws.Range("D4").Value = 3614531
$ws.Range("E4").Value = 2567912
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 191761

# --- Row 10: Sudafrica -------------------------------------------------------
$ws.Range("B10").Value = 635078
$ws.Range("C10").Value = 2063
$ws.Range("D10").Value = 557818
$ws.Range("E10").Value = 62582
$ws.Range("G10").Value = 115
$ws.Range("H10").Value = 14678

# --- Row 24: Alemania --------------------------------------------------------
$ws.Range("B24").Value = 250281
$ws.Range("C24").Value = 1467
$ws.Range("E24").Value = 16280
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 9401

# --- Row 54: Nigeria ---------------------------------------------------------
$ws.Range("B54").Value = 54587
$ws.Range("E54").Value = 10912

# --- Row 106: Luxemburgo ------------------------------------------------------
$ws.Range("B106").Value = 6854
$ws.Range("C106").Value = 43

# --- Row 107: Zimbabue --------------------------------------------------------
$ws.Range("B107").Value = 6837
$ws.Range("C107").Value = 159
$ws.Range("D107").Value = 5345
$ws.Range("E107").Value = 1286

# --- Rows 134/135: Angola overtakes Mali in rank -----------------------------
$ws.Range("A134").Value = "Angola"
$ws.Range("B134").Value = 2876
$ws.Range("C134").Value = 71
$ws.Range("D134").Value = 1167
$ws.Range("E134").Value = 1594
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = 115

$ws.Range("A135").Value = "Mali"
$ws.Range("B135").Value = 2807
$ws.Range("D135").Value = 2203
$ws.Range("E135").Value = 478
$ws.Range("H135").Value = 126

# --- Row 166: Republica del Chad ---------------------------------------------
$ws.Range("B166").Value = 1023
$ws.Range("C166").Value = 5
$ws.Range("D166").Value = 914
$ws.Range("E166").Value = 32

# --- Rows 214/215: Islas Malvinas overtakes Montserrat in rank ---------------
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
